# remove quit to menu button on pause when in menu and add ability to replay tutorial
#
# Spreadsheet-side of the change: the "StringLocalizations_BasicText" sheet
# gets a new localization row inserted (key + English text for the new
# "Replay Tutorial" button), and the active/selected tab moves from the
# Nicosia sheet to the BasicText sheet.

$wb = $excel.ActiveWorkbook

$basicText = $wb.Worksheets.Item("StringLocalizations_BasicText")

# Insert a new row above row 98 (inherits formatting from row 97, matching
# the existing localization rows) and fill in the new string pair.
$basicText.Rows.Item(98).Insert()

$basicText.Cells.Item(98, 1).Value = "BASIC_TEXT_REPLAY_TUTORIAL"
$basicText.Cells.Item(98, 2).Value = "Replay Tutorial"
$basicText.Cells.Item(98, 3).Value = "XXXX"
$basicText.Cells.Item(98, 4).Value = "XXXX"
$basicText.Cells.Item(98, 5).Value = "XXXX"

# Update the selection on this sheet to sit at the newly-entered row.
$basicText.Range("D99").Select()

# Make BasicText the active/selected tab (was Nicosia before).
$basicText.Activate()
